$d = $word.ActiveDocument

$replacements = @(
    @("35÷5=", "79÷8="),
    @("45÷4=", "49÷3="),
    @("81÷7=", "86÷9="),
    @("98÷9=", "20÷8="),
    @("72÷9=", "52÷7="),
    @("84÷7=", "65÷9="),
    @("59÷6=", "41÷6="),
    @("80÷8=", "66÷4="),
    @("89÷3=", "92÷3="),
    @("44÷3=", "33÷7="),
    @("46÷2=", "75÷9="),
    @("77÷5=", "68÷4="),
    @("21÷8=", "27÷7="),
    @("21÷2=", "15÷6="),
    @("78÷3=", "86÷2="),
    @("23÷2=", "82÷3="),
    @("94÷4=", "65÷2="),
    @("82÷4=", "90÷5="),
    @("91÷8=", "29÷6="),
    @("83÷3=", "34÷2="),
    @("29÷8=", "74÷3="),
    @("47÷2=", "80÷2="),
    @("12÷6=", "13÷5="),
    @("56÷9=", "26÷7="),
    @("54÷9=", "96÷7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
